# Updates cryptos list values per the Sat Sep  7 19:09:32 UTC 2024 GitHub Actions refresh.
# Set-TextValue writes a leading quote-prefix so Excel keeps Price-column cells as text
# instead of silently parsing them as numbers (which would drop formatting like trailing zeros).
function Set-TextValue($range, [string]$text) {
    $range.Value = '''' + $text
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '54.130.50'
$ws.Range('E2').Value = '  +0.91%  '
# Row 3
$ws.Range('D3').Value = '2.282.05'
$ws.Range('E3').Value = '  +2.60%  '
# Row 4
Set-TextValue $ws.Range('D4') '0.999'
$ws.Range('E4').Value = '  -0.04%  '
# Row 5
Set-TextValue $ws.Range('D5') '495.10'
$ws.Range('E5').Value = '  +2.58%  '
# Row 6
Set-TextValue $ws.Range('D6') '127.66'
$ws.Range('E6').Value = '  +1.52%  '
# Row 7
Set-TextValue $ws.Range('D7') '0.999'
$ws.Range('E7').Value = '  +0.02%  '
# Row 8
$ws.Range('E8').Value = '  +2.42%  '
# Row 9
$ws.Range('D9').Value = '2.279.42'
# Row 10
Set-TextValue $ws.Range('D10') '0.0948'
$ws.Range('E10').Value = '  +3.46%  '
# Row 11
$ws.Range('E11').Value = '  +2.33%  '
# Row 12
Set-TextValue $ws.Range('D12') '0.326'
$ws.Range('E12').Value = '  +3.64%  '
# Row 13
Set-TextValue $ws.Range('D13') '4.63'
$ws.Range('E13').Value = '  -0.97%  '
# Row 14
$ws.Range('D14').Value = '2.683.61'
$ws.Range('E14').Value = '  +2.48%  '
# Row 15
Set-TextValue $ws.Range('D15') '21.82'
$ws.Range('E15').Value = '  +3.64%  '
# Row 16
$ws.Range('D16').Value = '54.121.25'
$ws.Range('E16').Value = '  +1.06%  '
# Row 17
$ws.Range('E17').Value = '  +1.65%  '
# Row 18
$ws.Range('D18').Value = '2.289.25'
$ws.Range('E18').Value = '  +3.55%  '
# Row 19
$ws.Range('E19').Value = '  +5.15%  '
# Row 20
$ws.Range('E20').Value = '  +3.58%  '
# Row 21
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue $ws.Range('D21') '6.44'
$ws.Range('E21').Value = '  +5.66%  '
# Row 22
$ws.Range('B22').Value = 'BitcoinCash'
$ws.Range('C22').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue $ws.Range('D22') '300.62'
$ws.Range('E22').Value = '  +0.96%  '
# Row 23
Set-TextValue $ws.Range('D23') '0.998'
$ws.Range('E23').Value = '  -0.38%  '
# Row 24
Set-TextValue $ws.Range('D24') '5.38'
$ws.Range('E24').Value = '  -2.17%  '
# Row 25
Set-TextValue $ws.Range('D25') '62.48'
$ws.Range('E25').Value = '  -1.15%  '
# Row 26
Set-TextValue $ws.Range('D26') '1.02'
$ws.Range('E26').Value = '  +2.08%  '
# Row 27
$ws.Range('E27').Value = '  +2.63%  '
# Row 28
$ws.Range('D28').Value = '2.368.62'
$ws.Range('E28').Value = '  +1.68%  '
# Row 29
Set-TextValue $ws.Range('D29') '0.147'
$ws.Range('E29').Value = '  +3.41%  '
# Row 30
Set-TextValue $ws.Range('D30') '7.06'
$ws.Range('E30').Value = '  +1.35%  '
# Row 31
Set-TextValue $ws.Range('D31') '169.21'
$ws.Range('E31').Value = '  -0.09%  '
# Row 32
$ws.Range('E32').Value = '  +1.95%  '
# Row 33
$ws.Range('E33').Value = '  +1.16%  '
# Row 34
$ws.Range('E34').Value = '  +2.34%  '
# Row 35
Set-TextValue $ws.Range('D35') '0.999'
$ws.Range('E35').Value = '  +0.06%  '
# Row 36
Set-TextValue $ws.Range('D36') '0.997'
$ws.Range('E36').Value = '  +0.10%  '
# Row 37
$ws.Range('E37').Value = '  +1.96%  '
# Row 38
Set-TextValue $ws.Range('D38') '17.63'
$ws.Range('E38').Value = '  +1.50%  '
# Row 39
$ws.Range('E39').Value = '  +4.25%  '
# Row 40
Set-TextValue $ws.Range('D40') '0.861'
$ws.Range('E40').Value = '  +3.39%  '
# Row 41
$ws.Range('E41').Value = '  +3.55%  '
# Row 42
Set-TextValue $ws.Range('D42') '35.43'
$ws.Range('E42').Value = '  -0.67%  '
# Row 43
$ws.Range('E43').Value = '  +3.20%  '
# Row 44
$ws.Range('E44').Value = '  +2.75%  '
# Row 45
$ws.Range('E45').Value = '  +2.40%  '
# Row 46
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws.Range('D46') '127.38'
$ws.Range('E46').Value = '  +3.87%  '
# Row 47
$ws.Range('B47').Value = 'RenderToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range('D47') '4.77'
$ws.Range('E47').Value = '  +3.40%  '
# Row 48
Set-TextValue $ws.Range('D48') '0.0887'
$ws.Range('E48').Value = '  +1.35%  '
# Row 49
Set-TextValue $ws.Range('D49') '0.542'
$ws.Range('E49').Value = '  +1.72%  '
# Row 50
Set-TextValue $ws.Range('D50') '237.14'
$ws.Range('E50').Value = '  +3.57%  '
# Row 51
Set-TextValue $ws.Range('D51') '0.0483'
$ws.Range('E51').Value = '  +3.04%  '
